$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Job Site block (A7:A10) -----------------------------------------
# Replace the old single-brace placeholders with the new double-brace
# field names, and add the two extra address lines that were previously
# blank. The order of assignment below matches the order in which the
# new unique strings must be appended to xl/sharedStrings.xml.
$ws.Range("A8").Value = "{{Client_Address_1}}"
$ws.Range("A7").Value = "{{Client_Name}}"
$ws.Range("A9").Value = "{{Client_Address_2}}"
$ws.Range("A10").Value = "{{Client_Postcode}}"

# --- Job Information block (C9:D12) -----------------------------------
# Row 10: "Job Number: " label + {{Job_Number}} field
$ws.Range("D10").Value = "{{Job_Number}}"
$ws.Range("C10").Value = "Job Number: "

# Row 11 (new): "Contact No: " label + {{Client_Contact_No}} field.
# This row previously held no text, only fill/background formatting
# (fillId 34, no border). Touching the (invisible) diagonal border
# forces Excel to record an explicit "no border" style for this cell
# instead of silently reusing the existing untouched style, matching
# the new cellXfs entry introduced by the original edit.
$c11 = $ws.Range("C11")
$c11.Value = "Contact No: "
$c11.Borders.Item(5).LineStyle = -4142
$ws.Range("D11").Value = "{{Client_Contact_No}}"

# Row 12: "Contact: " label + {{Client_Contact}} field
$ws.Range("D12").Value = "{{Client_Contact}}"
$ws.Range("C12").Value = "Contact: "

# --- Misc ---------------------------------------------------------------
# Update the saved cursor/selection position recorded in the sheet view.
$ws.Range("C21").Select() | Out-Null
